# KSA_Cities.xlsx update
# 1. Fix the Arabic spelling of "Al Artawiah" (row 171, column C):
#    was "الأرطاية" -> now "الأرطاوية". Writing a plain literal also drops
#    the old (now-orphaned) shared string and the cell loses its border
#    style (s="1" -> no style), matching the authored edit.
# 2. Append three new city rows (177-179) after the last existing row (176):
#       Khbash    / خباش     | 17.557409          , 44.750014999999998  | Najran region  | South KSA
#       Qusaiba   / قصيباء   | 26.885753999999999 , 43.634118999999998  | Qassim region  | Middle KSA
#       Al Hmanah / الحمنة   | 23.008141999999999 , 39.883386000000002  | Madinah region | West KSA
#    The Area/Region text is copied (via .Text, verbatim, including the
#    non-breaking space inside the Qassim label) from existing rows that
#    already carry the exact same Arabic strings, so the new cells reuse
#    the existing shared-string entries instead of creating near-duplicates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct the Arabic name for "Al Artawiah" -----------------------
$ws.Range("C171").ClearFormats()
$ws.Range("C171").Value = "الأرطاوية"

# --- Reference cells already holding the exact Area / Region strings ----
$areaNajran  = $ws.Cells.Item(69, 6).Text   # F69  -> "منطقة نجران"
$areaQassim  = $ws.Cells.Item(31, 6).Text   # F31  -> "منطقة القصيم" (has NBSP)
$areaMadinah = $ws.Cells.Item(39, 6).Text   # F39  -> "منطقة المدينة المنورة"
$regionSouth = $ws.Cells.Item(10, 7).Text   # G10  -> "جنوب المملكة"
$regionMid   = $ws.Cells.Item(20, 7).Text   # G20  -> "وسط المملكة"
$regionWest  = $ws.Cells.Item(39, 7).Text   # G39  -> "غرب المملكة"

# --- 2. Add the new rows --------------------------------------------------
$newRows = @(
    @{ Row = 177; En = "Khbash";    Ar = "خباش";   Lat = 17.557409;           Lon = 44.750014999999998; Area = $areaNajran;  Region = $regionSouth },
    @{ Row = 178; En = "Qusaiba";   Ar = "قصيباء"; Lat = 26.885753999999999; Lon = 43.634118999999998; Area = $areaQassim;  Region = $regionMid },
    @{ Row = 179; En = "Al Hmanah"; Ar = "الحمنة"; Lat = 23.008141999999999; Lon = 39.883386000000002; Area = $areaMadinah; Region = $regionWest }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.En
    $ws.Cells.Item($row, 2).Value = $r.En
    $ws.Cells.Item($row, 3).Value = $r.Ar
    $ws.Cells.Item($row, 4).Value = $r.Lat
    $ws.Cells.Item($row, 5).Value = $r.Lon
    $ws.Cells.Item($row, 6).Value = $r.Area
    $ws.Cells.Item($row, 7).Value = $r.Region
}

# Copy the formatting (border style) of the last pre-existing data row onto
# the newly added rows so they match the rest of the table (style index 1).
$ws.Range("A176:G176").Copy()
$ws.Range("A177:G179").PasteSpecial(-4122)

# Keep the "select everything" selection in sync with the new used range.
$ws.Range("A1:G179").Select()
